$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column AU: header "11-ago" plus per-row values, continuing the
# existing daily-tracking table one more day.
$ws.Range("AU1").Value = "11-ago"

$values = @(17, 16, 12, 12, 13, 12, 13, 21, 17, 19)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 47).Value = $values[$i]
}

# Match the formatting used by the rest of the table (center-aligned
# integers for the data rows, text for the header).
$ws.Range("AU2:AU11").HorizontalAlignment = -4108
$ws.Range("AU2:AU11").NumberFormat = "0"
$ws.Range("AU1").NumberFormat = "@"

# Move selection to reflect post-edit active cell as seen in the diff
$ws.Range("AZ8").Select()
